# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Swap the two worker records (document number, name, and overdue amount)
# between row 16 and row 17 of the account-statement table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the two data rows.
$docNum16 = $ws.Range("C16").Value()
$name16   = $ws.Range("D16").Value()
$debt16   = $ws.Range("G16").Value()

$docNum17 = $ws.Range("C17").Value()
$name17   = $ws.Range("D17").Value()
$debt17   = $ws.Range("G17").Value()

# Write back the swapped data: row 16 now holds what used to be row 17
# (73164000 / MILTON JOSE ELITIM ZAMBRANO / 1423500) and row 17 now holds
# what used to be row 16 (33102898 / YENIS FONNEGRA BOLAÑOS / 1160000).
$ws.Range("C16").Value = $docNum17
$ws.Range("D16").Value = $name17
$ws.Range("G16").Value = $debt17

$ws.Range("C17").Value = $docNum16
$ws.Range("D17").Value = $name16
$ws.Range("G17").Value = $debt16
